$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: Vertex Cover Problem
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "Vertex_Cover_Problem"
$ws.Cells.Item(14, 3).Value = "DP"
$ws.Cells.Item(14, 4).Value = "Tree"
$ws.Cells.Item(14, 5).Value = "easy"
$ws.Cells.Item(14, 6).Value = "GeeksForGeeks"

# Row 15: Weighted Job Scheduling
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "Weighted Job Scheduling"
$ws.Cells.Item(15, 3).Value = "DP"
$ws.Cells.Item(15, 4).Value = "Array"
$ws.Cells.Item(15, 5).Value = "medium"
$ws.Cells.Item(15, 6).Value = "GeeksForGeeks"

# Apply the same centered style (style index 1) to the new cells, matching existing rows
$ws.Range("A14:F15").HorizontalAlignment = -4108  # xlCenter

# Update selection to match the recorded state after edit
$ws.Range("E22").Select()
